# Updates the line-power-flow results table (pl_mw) on rows 2-25 with the
# recomputed values from the "case with 380 kV done" run. Only columns
# B, C, E, F, G, I, J, K change; A, D, H, L, M, N, O stay as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers being updated, in order: B, C, E, F, G, I, J, K
# (columns A, D, H, L, M, N, O are unchanged)
$colNums = @(2, 3, 5, 6, 7, 9, 10, 11)

$data = @(
    @(0.4664733213809029, 0.1613733740226451, 0.2632282980981913, 3.332242289215031, 0.002530114545696351, 1.435168256707577, 0.1486094087404908, 0.623371105489241),
    @(0.4364923331386592, 0.1515352513802952, 0.2540395114709852, 3.272966288887346, 0.002534241460077589, 1.420343033369939, 0.1467668614985769, 0.5839119523902241),
    @(0.418386975326257, 0.145595935172679, 0.2486065196619478, 3.238296300987628, 0.00253690647417132, 1.4117461633018, 0.1456944492769878, 0.5600852147518651),
    @(0.4110850778917552, 0.1432008934074531, 0.2464450019855136, 3.224600290351589, 0.002538025558632251, 1.408369321945777, 0.1452721564031734, 0.5504762832931931),
    @(0.4098772039451148, 0.1428047193951869, 0.2460892506744301, 3.222352139521689, 0.00253821338250737, 1.407816215414385, 0.1452029215921584, 0.548886799595266),
    @(0.418288190958549, 0.1455635327192368, 0.2485771563208985, 3.238109843673271, 0.002536921432320737, 1.411700111044297, 0.1456886945919749, 0.5599552181244007),
    @(0.4560729379312818, 0.1579600321856276, 0.2600165846105611, 3.311444726877312, 0.00253151036839494, 1.42995106524809, 0.1479618136817749, 0.6096821156739622),
    @(0.5325799802582765, 0.18308415500519, 0.2841125425639959, 3.469032984882006, 0.002521934229289004, 1.469791599789431, 0.1528913540044883, 0.7103990764213677),
    @(0.5902756757171232, 0.2020562747520671, 0.3028398955767599, 3.593348484988127, 0.002515522460799113, 1.501586024072878, 0.1568074444625296, 0.7863825934459783),
    @(0.6168496362557221, 0.2108024600332215, 0.3115842618034037, 3.651786706073892, 0.002512739520663455, 1.516610030021511, 0.1586543529820261, 0.8213890046983465),
    @(0.6269598520048874, 0.2141313145436925, 0.3149280666228975, 3.674189189652424, 0.002511704819077146, 1.522380727129189, 0.1593632513142254, 0.8347089474817722),
    @(0.6247803340064593, 0.2134136324308429, 0.3142064708620893, 3.669352232533384, 0.002511926811018703, 1.521134269392533, 0.1592101526140368, 0.8318374175204326),
    @(0.6176804626526575, 0.2110759875869235, 0.3118587065170431, 3.653624285485165, 0.002512654012026231, 1.517083151997412, 0.158712483215119, 0.8224835642714083),
    @(0.6133377402631481, 0.2096463157025426, 0.3104248696012206, 3.644026099700199, 0.00251310193308856, 1.514612354976236, 0.1584088883743462, 0.8167623747414154),
    @(0.5885456098644966, 0.2014870382213303, 0.302272972341477, 3.589567534489817, 0.002515707016294048, 1.500615520280064, 0.1566880704952922, 0.7841037476143526),
    @(0.5734204790359172, 0.1965113751757599, 0.297329809364598, 3.556643502128395, 0.002517339347385274, 1.492173074590227, 0.1556492420196278, 0.7641819596519213),
    @(0.5647517617370283, 0.1936603790897777, 0.2945078220131165, 3.53788391712493, 0.00251829082161776, 1.487369916582111, 0.1550578906228708, 0.7527649740262916),
    @(0.561821980542959, 0.1926969414065809, 0.2935559813235642, 3.531562662512442, 0.002518615141616238, 1.485752678325426, 0.154858723298446, 0.7489064963361614),
    @(0.5750273807967403, 0.1970399159563669, 0.2978538236502146, 3.560129942615163, 0.002517164279440665, 1.493066325764389, 0.1557591893044474, 0.766298370481735),
    @(0.6197645843618318, 0.2117621510069796, 0.3125474187022803, 3.65823653482056, 0.002512439896845411, 1.51827084675395, 0.1588584017607744, 0.8252292841059159),
    @(0.6492782919654019, 0.2214823881169536, 0.3223400459923624, 3.72394817933349, 0.002509463740650441, 1.535218491101404, 0.1609394193676437, 0.864115937400129),
    @(0.6335010456939187, 0.216285433521989, 0.3170961546404953, 3.688730207218072, 0.002511042002851171, 1.526129475970507, 0.159823629006226, 0.8433272542891643),
    @(0.574300816659246, 0.1968009328509766, 0.2976168548076572, 3.558553196967353, 0.002517243387161994, 1.492662329612543, 0.1557094638153629, 0.7653414280549384),
    @(0.5116230377313684, 0.1761984771389677, 0.27741486535745, 3.424912149074402, 0.002524414767552026, 1.458574441513704, 0.1515065714058892, 0.6828060998317937)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $colNums.Count; $j++) {
        $ws.Cells.Item($row, $colNums[$j]).Value = $rowVals[$j]
    }
}

Write-Output "Updated $($data.Count) rows across $($colNums.Count) columns"
